# Update cryptos list (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sub5 = [char]0x2085
$sub8 = [char]0x2088

# Leading apostrophe forces Excel to treat the numeric-looking price
# strings as text, preserving exact digits (e.g. trailing zeros,
# multiple '.' separators) instead of letting them be parsed as numbers.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'26.107.35"
$ws.Range("E2").Value = "  -0.25%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.662.82"
$ws.Range("E3").Value = "  -0.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'207.97"
$ws.Range("E5").Value = "  -1.25%  "

# Row 6 - XRP
$ws.Range("D6").Value = "'0.5205"
$ws.Range("E6").Value = "  -1.49%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.2589"
$ws.Range("E8").Value = "  -3.41%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "'0.06313"
$ws.Range("E9").Value = "  +0.43%  "

# Row 10 - Solana
$ws.Range("D10").Value = "'20.97"
$ws.Range("E10").Value = "  -1.54%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.07534"
$ws.Range("E11").Value = "  +0.22%  "

# Row 12 - WrappedEther
$ws.Range("D12").Value = "'1.665.20"
$ws.Range("E12").Value = "  -0.73%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'4.410"
$ws.Range("E13").Value = "  -1.57%  "

# Row 14 - Polygon
$ws.Range("D14").Value = "'0.5363"
$ws.Range("E14").Value = "  -5.11%  "

# Row 15 & 16 - ShibaInu / Litecoin swap places
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'66.04"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0{0}7929" -f $sub5
$ws.Range("E16").Value = "  -2.17%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "'26.141.73"
$ws.Range("E17").Value = "  -0.24%  "

# Row 18 - Dai
$ws.Range("E18").Value = "  -0.19%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "'4.698"
$ws.Range("E19").Value = "  -3.08%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "'187.49"
$ws.Range("E20").Value = "  -0.54%  "

# Row 21 - Avalanche
$ws.Range("D21").Value = "'10.17"
$ws.Range("E21").Value = "  -3.30%  "

# Row 22 - Chainlink
$ws.Range("D22").Value = "'6.191"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23 - BinanceUSD
$ws.Range("E23").Value = "  -0.18%  "

# Row 24 - Monero
$ws.Range("D24").Value = "'148.84"
$ws.Range("E24").Value = "  +0.79%  "

# Row 25 - Stellar
$ws.Range("D25").Value = "'0.1223"
$ws.Range("E25").Value = "  -3.13%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.393"
$ws.Range("E26").Value = "  -2.79%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'15.64"
$ws.Range("E27").Value = "  -1.36%  "

# Row 28 - Hedera
$ws.Range("D28").Value = "'0.06153"
$ws.Range("E28").Value = "  -4.87%  "

# Row 29 - Toncoin
$ws.Range("D29").Value = "'1.370"
$ws.Range("E29").Value = "  +2.03%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -1.31%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'3.463"
$ws.Range("E31").Value = "  -1.80%  "

# Row 32 - Filecoin
$ws.Range("D32").Value = "'3.394"
$ws.Range("E32").Value = "  -2.38%  "

# Row 33 - LidoDAOToken
$ws.Range("D33").Value = "'1.630"
$ws.Range("E33").Value = "  -1.13%  "

# Row 34 - ARBITRUM
$ws.Range("D34").Value = "'0.9900"
$ws.Range("E34").Value = "  -1.54%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.95%  "

# Row 36 - MXToken
$ws.Range("D36").Value = "'2.752"
$ws.Range("E36").Value = "  +1.37%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'0.5869"
$ws.Range("E37").Value = "  -3.57%  "

# Row 38 - Maker
$ws.Range("D38").Value = "'1.104.89"
$ws.Range("E38").Value = "  +0.51%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "'0.01592"
$ws.Range("E39").Value = "  -1.57%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "'5.990"
$ws.Range("E40").Value = "  -2.73%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").Value = "'0.8446"
$ws.Range("E41").Value = "  -2.30%  "

# Row 42 - PaxDollar
$ws.Range("D42").Value = "'1.003"
$ws.Range("E42").Value = "  -0.30%  "

# Row 43 - Quant
$ws.Range("D43").Value = "'99.81"
$ws.Range("E43").Value = "  -0.18%  "

# Row 44 - RocketPoolETH
$ws.Range("D44").Value = "'1.812.47"
$ws.Range("E44").Value = "  -0.88%  "

# Row 45 - BabyDogeCoin
$ws.Range("D45").Value = "'0.0{0}108" -f $sub8
$ws.Range("E45").Value = "  -2.03%  "

# Row 46 - Aave
$ws.Range("D46").Value = "'55.15"

# Row 47 - Frax
$ws.Range("D47").Value = "'0.9997"
$ws.Range("E47").Value = "  -0.67%  "

# Row 48 - EnergySwap
$ws.Range("D48").Value = "'7.997"
$ws.Range("E48").Value = "  +0.49%  "

# Row 49 - Cronos
$ws.Range("D49").Value = "'0.05244"
$ws.Range("E49").Value = "  -0.43%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.47%  "

# Row 51 - Aptos
$ws.Range("D51").Value = "'5.865"
$ws.Range("E51").Value = "  -1.41%  "
